$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 22.52
$ws.Range("F2").Value = 9.82
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 76
$ws.Range("K2").Value = 63.8
$ws.Range("N2").Value = 85.96878041621773

# Row 3
$ws.Range("D3").Value = 307.13
$ws.Range("F3").Value = 5.76
$ws.Range("I3").Value = 63
$ws.Range("K3").Value = 63
$ws.Range("N3").Value = 85.96878041621773

# Row 4
$ws.Range("D4").Value = 25.57
$ws.Range("F4").Value = 12.15
$ws.Range("K4").Value = 59.8
$ws.Range("N4").Value = 85.96878041621773

# Row 5
$ws.Range("D5").Value = 48.9
$ws.Range("F5").Value = 19.27
$ws.Range("K5").Value = 57.2
$ws.Range("N5").Value = 85.96878041621773
